# Ajout d'un systeme d'authentification : remplace la ligne de campagne
# existante par un nouveau jeu de donnees (Produit A/B/C) et efface les
# lignes de campagnes obsoletes (3 a 6).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Ligne 2 : nouvelle campagne "C1IDacb615" ---
$ws.Range("A2").Value = 'C1IDacb615'
$ws.Range("B2").Value = '1:1'', 2:2'', 2:3'', 3:3'', 3:2'', 4:2'', 4:3'', 4:4'', 3:4'', 4:1'', 3:1'', 2:3'', 2:4'', 1:4'''
$ws.Range("C2").Value = 'produit B'', produit C'', produit C'', produit C'', produit C'', produit C'', produit C'', Produit A'', Produit A'', produit B'', produit B'', produit C'', Produit A'', Produit A'''
$ws.Range("D2").Value = 'u"Pr\xe9sence d''un label", produit durable'', produit durable'', Produit \xe9quitable'', Produit \xe9quitable'', Produit bio'', Produit bio'', Produit bio'', Produit \xe9quitable'', Produit bio'', Produit \xe9quitable'', produit durable'', produit durable'', u"Pr\xe9sence d''un label"'
$ws.Range("G2").Value = 'Produit A'

$ws.Range("E2").Value = 560.8571428571429
$ws.Range("F2").Value = 7852
$ws.Range("H2").Value = 1

# --- Lignes 3 a 6 : anciennes campagnes retirees (retour au token vide) ---
$ws.Range("A3:H6").Value = "Token"

